$d = $word.ActiveDocument

# Explicitly set PageBreakBefore = False (w:val="0") on every paragraph in
# the document body. Word's COM model normally omits pageBreakBefore when
# it is false; writing the property forces the explicit element/attribute
# to be (re)serialized, matching the "Keep formatting" toggle recorded in
# the commit.
foreach ($p in $d.Paragraphs) {
    $p.Format.PageBreakBefore = 0
}

# Also stamp the same explicit "no page break before" setting onto the
# built-in heading/title paragraph styles used by the document (these
# already carry keepNext/keepLines/spacing formatting, so Word's style
# writer adds pageBreakBefore alongside them rather than introducing a
# brand-new <w:pPr> the way it would for Normal/Table Normal).
$styleNames = @("Heading 1", "Heading 2", "Heading 3", "Heading 4", "Heading 5", "Heading 6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $s.ParagraphFormat.PageBreakBefore = 0
}
